$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Concepts")

# Version update
$ws1.Range("B3").Value = "6.0.0"

# Date update
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher row gets a value; remove duplicate Contact row, replace with Jurisdiction row
$ws1.Range("B9").Value = "Alvearie Team"
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

# Delete the now-duplicate row (old row 11 "Contact"/"No display for ContactDetail")
$ws1.Rows.Item(11).Delete()

# Case Sensitive gets value "true" (row 14 after the deletion shift)
# Plain .Value = "true" is auto-coerced to a boolean by Excel; route it
# through a formula + paste-special-values so it lands as literal text
# (matches the original cell's style, unlike the quote-prefix trick).
$ws1.Range("B14").Formula = "=""true"""
$ws1.Range("B14").Copy()
$ws1.Range("B14").PasteSpecial(-4163)
